$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 623
$ws.Range("F3").Value = 5885
$ws.Range("F4").Value = 72
$ws.Range("F7").Value = 1018
$ws.Range("F8").Value = 397
$ws.Range("F9").Value = 1375
$ws.Range("F11").Value = 3123
$ws.Range("F12").Value = 1959
$ws.Range("F13").Value = 122
$ws.Range("F15").Value = 201
$ws.Range("F16").Value = 82
$ws.Range("F17").Value = 175
$ws.Range("F19").Value = 994
$ws.Range("F20").Value = 364
$ws.Range("F22").Value = 63
$ws.Range("F23").Value = 3676
$ws.Range("F24").Value = 1167
$ws.Range("F25").Value = 2919
$ws.Range("F26").Value = 287
$ws.Range("F27").Value = 2273
$ws.Range("F28").Value = 4233
$ws.Range("F29").Value = 111
$ws.Range("F30").Value = 931
$ws.Range("F31").Value = 478
$ws.Range("F32").Value = 1342
$ws.Range("F33").Value = 97
$ws.Range("F34").Value = 6
$ws.Range("F35").Value = 28
$ws.Range("F36").Value = 26
$ws.Range("F37").Value = 35
$ws.Range("F38").Value = 1027
$ws.Range("F39").Value = 1290
$ws.Range("F41").Value = 1103
$ws.Range("F42").Value = 713
$ws.Range("F43").Value = 606
$ws.Range("F44").Value = 433
$ws.Range("F45").Value = 16
$ws.Range("F46").Value = 104
$ws.Range("F47").Value = 3
$ws.Range("F48").Value = 323
$ws.Range("F49").Value = 3611

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 917
$ws.Range("F25").Value = 18
$ws.Range("F28").Value = 44

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 623
$ws.Range("F3").Value = 5886
$ws.Range("F4").Value = 72
$ws.Range("F7").Value = 397
$ws.Range("F8").Value = 1375
$ws.Range("F9").Value = 3123
$ws.Range("F11").Value = 1959
$ws.Range("F12").Value = 122
$ws.Range("F14").Value = 201
$ws.Range("F15").Value = 917
$ws.Range("F17").Value = 82
$ws.Range("F18").Value = 175
$ws.Range("F19").Value = 994
$ws.Range("F20").Value = 364
$ws.Range("F21").Value = 3676
$ws.Range("F23").Value = 1167
$ws.Range("F25").Value = 2919
$ws.Range("F26").Value = 2273
$ws.Range("F27").Value = 4233
$ws.Range("F28").Value = 111
$ws.Range("F29").Value = 931
$ws.Range("F30").Value = 1342
$ws.Range("F31").Value = 35
$ws.Range("F32").Value = 1027
$ws.Range("F34").Value = 1290
$ws.Range("F36").Value = 1103
$ws.Range("F38").Value = 713
$ws.Range("F40").Value = 433
$ws.Range("F43").Value = 16
$ws.Range("F44").Value = 18
$ws.Range("F45").Value = 104
$ws.Range("F47").Value = 323
$ws.Range("F48").Value = 3611
$ws.Range("F49").Value = 44
